$d = $word.ActiveDocument

$replacements = @(
    @{old="71÷2="; new="29÷9="},
    @{old="68÷8="; new="45÷7="},
    @{old="49÷7="; new="86÷8="},
    @{old="96÷4="; new="71÷8="},
    @{old="86÷2="; new="75÷6="},
    @{old="12÷3="; new="73÷5="},
    @{old="60÷9="; new="77÷5="},
    @{old="50÷9="; new="67÷8="},
    @{old="51÷7="; new="12÷8="},
    @{old="59÷6="; new="59÷9="},
    @{old="80÷2="; new="33÷8="},
    @{old="32÷5="; new="10÷6="},
    @{old="71÷9="; new="76÷2="},
    @{old="65÷7="; new="42÷6="},
    @{old="37÷4="; new="54÷9="},
    @{old="68÷4="; new="39÷4="},
    @{old="30÷4="; new="86÷6="},
    @{old="43÷4="; new="43÷6="},
    @{old="80÷9="; new="92÷8="},
    @{old="21÷7="; new="49÷3="},
    @{old="65÷6="; new="83÷9="},
    @{old="78÷7="; new="11÷7="},
    @{old="82÷9="; new="86÷5="},
    @{old="17÷2="; new="16÷4="},
    @{old="15÷6="; new="30÷9="}
)

foreach ($r in $replacements) {
    $d.Content.Find.Execute($r.old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $r.new, 2)
}
